$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append at row 54 (2019-04-14)
$rowIndex = 54
$values = @(43569, 1, 1, 5, 11, 3, 1, 1, 440, 1, 1, 1, 2, 1, 2, 1, 0, 2, 2, 1)

# Column A holds a date serial value, formatted the same way as the cells above it.
# Copy A53's format (style index 1, numFmtId 14 = m/d/yyyy) onto A54 via PasteSpecial
# so we reuse the existing style entry instead of minting a new numFmt/xf.
$dateCell = $ws.Cells.Item($rowIndex, 1)
$dateCell.Value = 43569
$ws.Cells.Item($rowIndex - 1, 1).Copy()
$dateCell.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($col = 2; $col -le 20; $col++) {
    $ws.Cells.Item($rowIndex, $col).Value = $values[$col - 1]
}

# Update the active selection to match the post-edit state: the whole row 54 is
# selected, the same way it was selected before the new data was typed in (only
# the active cell within that selection advances, from A54 to B54, once the row
# has been filled in).
$ws.Rows.Item(54).Select()

